# ieee9_ESS.xlsx — "ESS cost data. 2025 and 2035"
# Adds two new year columns (2025, 2035) to the "Investment Cost" sheet,
# interpolated as the average of their neighbouring (existing) years, and
# highlights those two new columns with a themed fill. Also nudges a
# couple of sheet selections / the active sheet, matching the authored
# commit.

$wb = $excel.ActiveWorkbook

$wsNREL = $wb.Worksheets.Item("Investment Cost NREL")
$wsCost = $wb.Worksheets.Item("Investment Cost")

# ---------------------------------------------------------------------
# 1) "Investment Cost" sheet: insert a 2025 column before the existing
#    2030 column, and a 2035 column before the existing 2040 column.
#    (Inserting shifts the existing B:E columns — with their formulas —
#    out to B,D,F,G automatically, preserving per-cell formatting.)
# ---------------------------------------------------------------------

# Insert column for 2025 (becomes column C; old C..E shift to D..F)
$wsCost.Columns.Item(3).Insert()
# Insert column for 2035 (becomes column E; old D..E [now at D..E] shift to E..F->F..G)
$wsCost.Columns.Item(5).Insert()

# Match the (best-fit) width of the surrounding data columns.
$wsCost.Columns.Item(3).ColumnWidth = 8.67
$wsCost.Columns.Item(5).ColumnWidth = 8.67

# Header years
$wsCost.Range("C1").Value = 2025
$wsCost.Range("E1").Value = 2035

# Interpolated cost rows: simple average of the two surrounding years
$wsCost.Range("C2").Formula = "=AVERAGE(D2,B2)"
$wsCost.Range("E2").Formula = "=AVERAGE(F2,D2)"
$wsCost.Range("C3").Formula = "=AVERAGE(D3,B3)"
$wsCost.Range("E3").Formula = "=AVERAGE(F3,D3)"

# Highlight the two new (interpolated) columns with a theme fill so
# they're visually distinguished from the sourced NREL data.
$wsCost.Range("C2").Interior.ThemeColor = 4
$wsCost.Range("C3").Interior.ThemeColor = 4
$wsCost.Range("E2").Interior.ThemeColor = 4
$wsCost.Range("E3").Interior.ThemeColor = 4

# ---------------------------------------------------------------------
# 2) Selections / active sheet bookkeeping
# ---------------------------------------------------------------------

# "Investment Cost NREL" sheet: move the lingering selection back onto
# the data table.
[void]$wsNREL.Range("D4").Select()

# "Investment Cost" sheet: select the newly-built data block and make
# this the active tab.
[void]$wsCost.Range("A2:G3").Select()
$wsCost.Activate()
